$d = $word.ActiveDocument

# 1) Date line: "Tlaxcoapan, Hgo" -> "Municipio, Estado"
#    and "de 2024" -> "de Año" (only the trailing "de <year>" after the {fechmes} placeholder)
$rng = $d.Content
$rng.Find.Execute("Tlaxcoapan, Hgo", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "Municipio, Estado", 2) | Out-Null

$rng = $d.Content
$rng.Find.Execute("de 2024", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "de Año", 2) | Out-Null

# 2) "es propiedad de la Escuela Preparatoria Número 6" -> "es propiedad de la Escuela Nombre de la Escuela"
$rng = $d.Content
$rng.Find.Execute("Escuela Preparatoria Número 6.", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "Escuela Nombre de la Escuela.", 2) | Out-Null

# 3) Remove the signature text box (shape) in the body
while ($d.Shapes.Count -gt 0) {
    $d.Shapes.Item(1).Delete()
}

# 4) Remove the background image from the default header (header2.xml)
$sec = $d.Sections.Item(1)
$hdr = $sec.Headers.Item(1)
while ($hdr.Shapes.Count -gt 0) {
    $hdr.Shapes.Item(1).Delete()
}
